# Updates cryptos list figures (Price / Volume(1h)) to the latest scrape.
# Cells are plain text (inline strings) in the workbook, so each write forces
# Text number-formatting before assignment (and restores the default "Normal"
# style afterwards) to stop Excel from silently re-interpreting decimal-looking
# price strings (e.g. "8.10", "1.00") as numbers, which would drop the text
# formatting/trailing zeros that the source data relies on.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "59.667.17"
Set-TextValue "E2" "  -2.31%  "
Set-TextValue "D3" "2.368.89"
Set-TextValue "E3" "  -2.69%  "
Set-TextValue "E4" "  +0.18%  "
Set-TextValue "D5" "559.23"
Set-TextValue "E5" "  -2.29%  "
Set-TextValue "D6" "137.74"
Set-TextValue "E6" "  -2.05%  "
Set-TextValue "D8" "0.529"
Set-TextValue "E8" "  -0.36%  "
Set-TextValue "D9" "2.365.86"
Set-TextValue "D10" "0.104"
Set-TextValue "E10" "  -3.95%  "
Set-TextValue "E11" "  -1.15%  "
Set-TextValue "E12" "  -1.47%  "
Set-TextValue "E13" "  -1.12%  "
Set-TextValue "E14" "  -2.37%  "
Set-TextValue "D15" "2.799.06"
Set-TextValue "E15" "  -2.01%  "
Set-TextValue "D16" "0.0000164"
Set-TextValue "E16" "  -3.58%  "
Set-TextValue "D17" "59.425.50"
Set-TextValue "E17" "  -2.65%  "
Set-TextValue "D18" "2.368.17"
Set-TextValue "E18" "  -3.08%  "
Set-TextValue "D19" "8.10"
Set-TextValue "E19" "  +11.06%  "
Set-TextValue "D20" "10.46"
Set-TextValue "E20" "  -1.34%  "
Set-TextValue "D21" "321.27"
Set-TextValue "E21" "  -0.90%  "
Set-TextValue "D22" "4.02"
Set-TextValue "E22" "  -0.62%  "
Set-TextValue "E23" "  -3.33%  "
Set-TextValue "E24" "  +0.04%  "
Set-TextValue "E25" "  -4.62%  "
Set-TextValue "D26" "64.14"
Set-TextValue "D27" "552.04"
Set-TextValue "E27" "  -5.06%  "
Set-TextValue "D28" "8.08"
Set-TextValue "E28" "  -8.97%  "
Set-TextValue "D29" "2.482.40"
Set-TextValue "E29" "  -2.98%  "
Set-TextValue "E30" "  -0.19%  "
Set-TextValue "E31" "  +1.10%  "
Set-TextValue "E32" "  -3.90%  "
Set-TextValue "E33" "  -3.83%  "
Set-TextValue "E34" "  -2.31%  "
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  -0.41%  "
Set-TextValue "E36" "  +1.81%  "
Set-TextValue "D37" "151.19"
Set-TextValue "E37" "  +0.31%  "
Set-TextValue "D38" "0.365"
Set-TextValue "E38" "  -1.06%  "
Set-TextValue "E39" "  -2.22%  "
Set-TextValue "E40" "  -1.02%  "
Set-TextValue "E41" "  -2.67%  "
Set-TextValue "E42" "  -0.06%  "
Set-TextValue "D43" "41.38"
Set-TextValue "E43" "  -0.79%  "
Set-TextValue "E44" "  -1.12%  "
Set-TextValue "E45" "  +1.39%  "
Set-TextValue "D46" "0.0₆0299"
Set-TextValue "E46" "  +4.81%  "
Set-TextValue "D47" "137.97"
Set-TextValue "E47" "  -2.66%  "
Set-TextValue "E48" "  -0.79%  "
Set-TextValue "E50" "  -1.66%  "
Set-TextValue "D51" "19.06"
Set-TextValue "E51" "  -2.92%  "
